# Generate Report for Handback
# Refreshes the "Latest HO Xliff Generate Date" / handoff / handback
# timestamp columns to reflect a fresh handback report run.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
# (24c6431a-...) moves from 07:10:23 to 07:11:10.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 07:11:10"

# zh-cn sheet: handoff/handback datetimes for the first file move from
# 07:10:18 / 07:10:35 to 07:11:01 / 07:11:59.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 07:11:01"
$wsZhCn.Range("K2").Value = "2016-09-04 07:11:59"

# de-de sheet: handoff datetime for the first file mirrors the Overview
# sheet's value (07:10:23 -> 07:11:10); handback datetime moves from
# 07:10:42 to 07:12:12.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 07:11:10"
$wsDeDe.Range("K2").Value = "2016-09-04 07:12:12"
